$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 12:35"

# Swap the province labels for rows 34 and 35 (Soria <-> Gran Canaria)
$ws.Range("A34").Value = "Gran Canaria"
$ws.Range("A35").Value = "Soria"

# Update numeric data cells (row => B,C,D,E values)
$ws.Range("B4").Value = 64972
$ws.Range("C4").Value = 39751
$ws.Range("D4").Value = 16501
$ws.Range("E4").Value = 8720

$ws.Range("B5").Value = 55144
$ws.Range("C5").Value = 24966
$ws.Range("D5").Value = 24557
$ws.Range("E5").Value = 5621

$ws.Range("B6").Value = 18048
$ws.Range("C6").Value = 7472
$ws.Range("D6").Value = 8662
$ws.Range("E6").Value = 1914

$ws.Range("B7").Value = 16431
$ws.Range("C7").Value = 6181
$ws.Range("D7").Value = 7439
$ws.Range("E7").Value = 2811

$ws.Range("B9").Value = 12293
$ws.Range("C9").Value = 9297
$ws.Range("D9").Value = 1670
$ws.Range("E9").Value = 1326

$ws.Range("B11").Value = 9291
$ws.Range("C11").Value = 7820
$ws.Range("D11").Value = 876
$ws.Range("E11").Value = 595

$ws.Range("B14").Value = 5359
$ws.Range("C14").Value = 3186
$ws.Range("D14").Value = 1345
$ws.Range("E14").Value = 828

$ws.Range("B16").Value = 5081
$ws.Range("C16").Value = 3155
$ws.Range("D16").Value = 1433
$ws.Range("E16").Value = 493

$ws.Range("C23").Value = 2749
$ws.Range("D23").Value = 912
$ws.Range("E23").Value = 347

$ws.Range("B29").Value = 2919
$ws.Range("C29").Value = 2422
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = 487

$ws.Range("B33").Value = 2351
$ws.Range("C33").Value = 1031
$ws.Range("D33").Value = 1016
$ws.Range("E33").Value = 304

$ws.Range("B34").Value = 2268
$ws.Range("C34").Value = 1456
$ws.Range("D34").Value = 661
$ws.Range("E34").Value = 151

$ws.Range("B35").Value = 2262
$ws.Range("C35").Value = 380
$ws.Range("D35").Value = 1764
$ws.Range("E35").Value = 118

$ws.Range("C36").Value = 1981
$ws.Range("D36").Value = 62
$ws.Range("E36").Value = 203

$ws.Range("C59").Value = 115
$ws.Range("D59").Value = 2
